$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Insert a brand new row at the very top for "Date and Time".
#    This shifts all the existing rows (1-42) down to (2-43).
# -----------------------------------------------------------------
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "Date and Time"
$ws.Range("B1").Value = "2024-03-12 14:51:07.889000 to 2024-03-12 15:56:34.474000"

# -----------------------------------------------------------------
# 2. Relabel / rename cells whose text changed (units appended etc.)
#    Row numbers below refer to the *current* sheet state (post top-insert).
# -----------------------------------------------------------------
$ws.Range("A9").Value  = "Total distance covered (km)"
$ws.Range("A10").Value = "Total energy consumption(WH/KM)"
$ws.Range("A11").Value = "Total SOC consumed(%)"
$ws.Range("A13").Value = "Peak Power(kW)"
$ws.Range("A14").Value = "Average Power(kW)"
$ws.Range("A15").Value = "Total Energy Regenerated(kWh)"
$ws.Range("A16").Value = "Regenerative Effectiveness(%)"
$ws.Range("B16").Value = 0.007660860966486528

# Highest / Lowest Cell Voltage swap order (and relabel with units)
$ws.Range("A17").Value = "Highest Cell Voltage(V)"
$ws.Range("B17").Value = 3.414
$ws.Range("A18").Value = "Lowest Cell Voltage(V)"
$ws.Range("B18").Value = 3.107

$ws.Range("A19").Value = "Difference in Cell Voltage(V)"
$ws.Range("A20").Value = "Minimum Temperature(C)"
$ws.Range("A21").Value = "Maximum Temperature(C)"
$ws.Range("A22").Value = "Difference in Temperature(C)"
$ws.Range("B22").Value = 10
$ws.Range("A23").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("A24").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("A25").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("A26").Value = "Maximum MCU Temperature(C)"
$ws.Range("A27").Value = "Maximum Motor Temperature(C)"
$ws.Range("A28").Value = "Abnormal Motor Temperature Detected(C)"

# highest / lowest cell temp swap order (and relabel with units)
$ws.Range("A29").Value = "highest cell temp(C)"
$ws.Range("B29").Value = 46
$ws.Range("A30").Value = "lowest cell temp(C)"
$ws.Range("B30").Value = -1

$ws.Range("A31").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"

# -----------------------------------------------------------------
# 3. Remove the now-redundant "Maximum BMS Temperature in C" row.
#    This shifts everything below it up by one.
# -----------------------------------------------------------------
$ws.Rows.Item(32).Delete()

# -----------------------------------------------------------------
# 4. Continue relabeling / updating values (rows shifted up by one now).
# -----------------------------------------------------------------
$ws.Range("A32").Value = "Battery Voltage(V)"
$ws.Range("B32").Value = 53

$ws.Range("A33").Value = "Total energy charged(kWh)"
$ws.Range("B33").Value = 1.560140433055555

$ws.Range("A34").Value = "Electricity consumption units(kW)"
$ws.Range("B34").Value = 0.000000110836916244356

# -----------------------------------------------------------------
# 5. Insert a new row for "Cycle Count of battery" right before
#    "Idling time percentage".
# -----------------------------------------------------------------
$ws.Rows.Item(35).Insert()
$ws.Range("A35").Value = "Cycle Count of battery"
$ws.Range("B35").Value = 54

# -----------------------------------------------------------------
# 6. Update the remaining "Idling time percentage" / "Time spent in ..."
#    values (now at rows 36-43).
# -----------------------------------------------------------------
$ws.Range("B36").Value = 8.756917108791146   # Idling time percentage

$ws.Range("B37").Value = 9.695361971589938   # Time spent in 0-10 km/h
$ws.Range("B38").Value = 3.465685435563923   # Time spent in 10-20 km/h
$ws.Range("B39").Value = 5.125791545438987   # Time spent in 20-30 km/h
$ws.Range("B40").Value = 21.672656740259     # Time spent in 30-40 km/h
$ws.Range("B41").Value = 10.04050430714815   # Time spent in 40-50 km/h
$ws.Range("B42").Value = 34.39157966797877   # Time spent in 50-60 km/h
$ws.Range("B43").Value = 6.768783159335958   # Time spent in 60-70 km/h

# -----------------------------------------------------------------
# 7. Append two brand new rows for the higher speed buckets.
# -----------------------------------------------------------------
$ws.Range("A44").Value = "Time spent in 70-80 km/h"
$ws.Range("B44").Value = 0
$ws.Range("A45").Value = "Time spent in 80-90 km/h"
$ws.Range("B45").Value = 0
